$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 10-19 (columns A-M) with the new averaged-intensity data
# reflecting the added spiral sampling schemes and the shifted
# Gaussian-Quadrature / NoRotation.. / Rotation.. / HexGrid.. rows.

# Row 10: Gaussian-Quadrature
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 'Gaussian-Quadrature'
$ws.Cells.Item(10,3).Value = 0.9917500818181104
$ws.Cells.Item(10,4).Value = 1.00398004565772
$ws.Cells.Item(10,5).Value = 0.9983221536969686
$ws.Cells.Item(10,6).Value = 0.9917500818181104
$ws.Cells.Item(10,7).Value = 0.9993463988198579
$ws.Cells.Item(10,8).Value = 0.9929919182923702
$ws.Cells.Item(10,9).Value = 0.9952941176470588
$ws.Cells.Item(10,10).Value = 1.00398004565772
$ws.Cells.Item(10,11).Value = 1.001151099677344
$ws.Cells.Item(10,12).Value = 0.9964505907477272
$ws.Cells.Item(10,13).Value = 0.9969474526553476

# Row 11: Spiral-90deg-10rot-5space
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 'Spiral-90deg-10rot-5space'
$ws.Cells.Item(11,3).Value = 0.9998580257616398
$ws.Cells.Item(11,4).Value = 0.9803478011591985
$ws.Cells.Item(11,5).Value = 0.9999787357187869
$ws.Cells.Item(11,6).Value = 0.9998580257616398
$ws.Cells.Item(11,7).Value = 0.9851376176462555
$ws.Cells.Item(11,8).Value = 1.010729084797579
$ws.Cells.Item(11,9).Value = 1.002094307701241
$ws.Cells.Item(11,10).Value = 0.9803478011591985
$ws.Cells.Item(11,11).Value = 0.9901632684389927
$ws.Cells.Item(11,12).Value = 0.9950106471003162
$ws.Cells.Item(11,13).Value = 0.9963575954641168

# Row 12: Spiral-90deg-15rot-5space
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 'Spiral-90deg-15rot-5space'
$ws.Cells.Item(12,3).Value = 0.9997889145754573
$ws.Cells.Item(12,4).Value = 0.9806081777437193
$ws.Cells.Item(12,5).Value = 0.9999233910826922
$ws.Cells.Item(12,6).Value = 0.9997889145754573
$ws.Cells.Item(12,7).Value = 0.9853164160587775
$ws.Cells.Item(12,8).Value = 1.010457457135877
$ws.Cells.Item(12,9).Value = 1.002040855315162
$ws.Cells.Item(12,10).Value = 0.9806081777437193
$ws.Cells.Item(12,11).Value = 0.9902657844132058
$ws.Cells.Item(12,12).Value = 0.9950273494943316
$ws.Cells.Item(12,13).Value = 0.9963558686519476

# Row 13: Spiral-90deg-10rot-3space
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 'Spiral-90deg-10rot-3space'
$ws.Cells.Item(13,3).Value = 0.9998549418076437
$ws.Cells.Item(13,4).Value = 0.9804208854006838
$ws.Cells.Item(13,5).Value = 0.9999661422583563
$ws.Cells.Item(13,6).Value = 0.9998549418076437
$ws.Cells.Item(13,7).Value = 0.9851805794357783
$ws.Cells.Item(13,8).Value = 1.01061473478673
$ws.Cells.Item(13,9).Value = 1.002077917222269
$ws.Cells.Item(13,10).Value = 0.9804208854006838
$ws.Cells.Item(13,11).Value = 0.9901935138295201
$ws.Cells.Item(13,12).Value = 0.9950242278185819
$ws.Cells.Item(13,13).Value = 0.9963525334852434

# Row 14: NoRotation-tilt60deg
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 'NoRotation-tilt60deg'
$ws.Cells.Item(14,3).Value = 0.9891759999999986
$ws.Cells.Item(14,4).Value = 1.024568000000002
$ws.Cells.Item(14,5).Value = 0.9856279999999995
$ws.Cells.Item(14,6).Value = 0.9891759999999986
$ws.Cells.Item(14,7).Value = 1.013736
$ws.Cells.Item(14,8).Value = 0.9939719999999989
$ws.Cells.Item(14,9).Value = 0.9916720000000001
$ws.Cells.Item(14,10).Value = 1.024568000000002
$ws.Cells.Item(14,11).Value = 1.005098000000001
$ws.Cells.Item(14,12).Value = 0.9971369999999997
$ws.Cells.Item(14,13).Value = 0.9997919999999998

# Row 15: Rotation-NoTilt
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 'Rotation-NoTilt'
$ws.Cells.Item(15,3).Value = 0.9936249999999993
$ws.Cells.Item(15,4).Value = 1.03
$ws.Cells.Item(15,5).Value = 0.98
$ws.Cells.Item(15,6).Value = 0.9936249999999993
$ws.Cells.Item(15,7).Value = 1.02
$ws.Cells.Item(15,8).Value = 1
$ws.Cells.Item(15,9).Value = 0.99
$ws.Cells.Item(15,10).Value = 1.03
$ws.Cells.Item(15,11).Value = 1.005
$ws.Cells.Item(15,12).Value = 0.9993124999999996
$ws.Cells.Item(15,13).Value = 1.002270833333333

# Row 16: Rotation-60detTilt
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 'Rotation-60detTilt'
$ws.Cells.Item(16,3).Value = 0.9952883510272009
$ws.Cells.Item(16,4).Value = 1.016697255935997
$ws.Cells.Item(16,5).Value = 0.9874073931776012
$ws.Cells.Item(16,6).Value = 0.9952883510272009
$ws.Cells.Item(16,7).Value = 1.010231527833597
$ws.Cells.Item(16,8).Value = 0.9975516719104014
$ws.Cells.Item(16,9).Value = 0.9929900587008014
$ws.Cells.Item(16,10).Value = 1.016697255935997
$ws.Cells.Item(16,11).Value = 1.002052324556799
$ws.Cells.Item(16,12).Value = 0.9986703377920001
$ws.Cells.Item(16,13).Value = 1.000027709764266

# Row 17: HexGrid-90degTilt5degRes
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 'HexGrid-90degTilt5degRes'
$ws.Cells.Item(17,3).Value = 0.9973810914737039
$ws.Cells.Item(17,4).Value = 0.9977039123396846
$ws.Cells.Item(17,5).Value = 0.9972887469786561
$ws.Cells.Item(17,6).Value = 0.9973810914737039
$ws.Cells.Item(17,7).Value = 0.9974735170355726
$ws.Cells.Item(17,8).Value = 0.997582915918265
$ws.Cells.Item(17,9).Value = 0.9975960388977326
$ws.Cells.Item(17,10).Value = 0.9977039123396846
$ws.Cells.Item(17,11).Value = 0.9974963296591703
$ws.Cells.Item(17,12).Value = 0.9974387105664372
$ws.Cells.Item(17,13).Value = 0.9975043704406025

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 'HexGrid-90degTilt22p5degRes'
$ws.Cells.Item(18,3).Value = 0.9996018759134899
$ws.Cells.Item(18,4).Value = 0.9957941903194829
$ws.Cells.Item(18,5).Value = 1.000558675344193
$ws.Cells.Item(18,6).Value = 0.9996018759134899
$ws.Cells.Item(18,7).Value = 0.9973942844441599
$ws.Cells.Item(18,8).Value = 0.9993231458200004
$ws.Cells.Item(18,9).Value = 0.9983616986462573
$ws.Cells.Item(18,10).Value = 0.9957941903194829
$ws.Cells.Item(18,11).Value = 0.998176432831838
$ws.Cells.Item(18,12).Value = 0.998889154372664
$ws.Cells.Item(18,13).Value = 0.9985056450812638

# Row 19: HexGrid-60degTilt5degRes
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 'HexGrid-60degTilt5degRes'
$ws.Cells.Item(19,3).Value = 0.9978881304532299
$ws.Cells.Item(19,4).Value = 0.9947627792378595
$ws.Cells.Item(19,5).Value = 0.9981874446985601
$ws.Cells.Item(19,6).Value = 0.9978881304532299
$ws.Cells.Item(19,7).Value = 0.9949480237959699
$ws.Cells.Item(19,8).Value = 1.000228881941848
$ws.Cells.Item(19,9).Value = 0.9981252081641926
$ws.Cells.Item(19,10).Value = 0.9947627792378595
$ws.Cells.Item(19,11).Value = 0.9964751119682098
$ws.Cells.Item(19,12).Value = 0.9971816212107198
$ws.Cells.Item(19,13).Value = 0.9973567447152766

# Column A on the newly-added rows needs the same style as the rest of column A
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = 0
